$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: Description in column A, Command in column B
$ws.Range("A6").Value = "Checkout remote branch"
$ws.Range("B6").Value = "git fetch`ngit checkout [name of remote branch]"

# Wrap text for the new command cell and size the row to fit
$ws.Range("B6").WrapText = $true
$ws.Range("A6:B6").RowHeight = 30

# Update the selection to reflect where the user ended up
$ws.Range("C6").Select() | Out-Null
